$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$dst = $ws.Range("E15")
$dst.Value = "***.*"
$src = $ws.Range("A14")
$src.Copy()
$dst.PasteSpecial(-4122)
